$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (2-7) for columns E,G,H,I,J,K,M,N,O,P,Q,R,S,T
$data = @{
    2 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989;
           K=3; M=6.878910333333334; N=20.636731; O=0.14811996585983; P=0.14811996585983;
           Q=62.98990217997979; R=566.9091196198181; S=0.1435977838749157; T=0.1435977838749156 }
    3 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989;
           K=3; M=22.22983366666666; N=66.68950099999999; O=0.4786633411720634; P=0.4786633411720634;
           Q=203.5576828724309; R=1832.019145851878; S=0.464049492689708; T=0.464049492689708 }
    4 = @{ E=3; G=9.156959333333335; H=27.470878; I=0.969469463764299; J=0.9694694637642989;
           K=3; M=17.33273533333334; N=51.998206; O=0.3732166929681066; P=0.3732166929681066;
           Q=158.7151525827632; R=1428.436373244868; S=0.3618221871996753; T=0.3618221871996753 }
    5 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109;
           K=3; M=6.878910333333334; N=20.636731; O=0.14811996585983; P=0.14811996585983;
           Q=1.983678251733667; R=17.853104265603; S=0.004522181984914348; T=0.004522181984914348 }
    6 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109;
           K=3; M=22.22983366666666; N=66.68950099999999; O=0.4786633411720634; P=0.4786633411720634;
           Q=6.410439364290332; R=57.69395427861299; S=0.01461384848235543; T=0.01461384848235544 }
    7 = @{ E=3; G=0.288371; H=0.865113; I=0.03053053623570109; J=0.03053053623570109;
           K=3; M=17.33273533333334; N=51.998206; O=0.3732166929681066; P=0.3732166929681066;
           Q=4.998258220808667; R=44.984323987278; S=0.01139450576843131; T=0.01139450576843131 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
